$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "317.00"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3.30%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "20"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.00"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.00%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "20"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.136"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.59%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "20"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08239"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.88%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "20"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.071"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "6.69%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "20"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.323"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.47%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "20"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.324"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "3.25%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "20"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9407"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.01%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "20"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1373"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-6.77%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "20"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1980"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.52%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "20"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09106"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.32%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "20"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03504"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.02%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "20"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09814"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.15%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "20"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001374"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.47%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "20"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006252"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "6.76%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "20"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.707"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.10%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "20"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.360"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.82%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "20"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3496"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.05%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "20"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1308"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.45%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "20"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.998"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "9.50%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "20"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2444"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.28%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "20"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04352"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.63%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "20"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001228"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.67%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "20"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004829"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "12.78%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "20"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001301"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.07%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "20"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003990"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-10.28%"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "20"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "20"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "20"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "20"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "20"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "20"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "20"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "20"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "20"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "20"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "20"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "20"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02220"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "8.71%"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "20"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05212"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.70%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "20"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007680"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.15%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "20"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009769"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-5.20%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "20"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1409"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.37%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "20"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002049"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-3.44%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "20"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008897"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.13%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "20"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006596"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "6.45%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "20"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "20"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.001686"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "5.27%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "20"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002935"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-5.62%"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "20"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "20"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "20"
